$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11: Marking right/wrong values corrected
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12: Total marks corrected
$ws.Range("B12").Value = 100
$ws.Range("E12").Value = "100 / 112"
